# Update benchmark: 2025-12-30 06:42:47 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ŞANS OYUNLARI
$ws.Range("C2").Value = "25 TL - 25 TL"
$ws.Range("F2").Value = "28,57 TL - 28,57 TL"

# Row 3 - HESAPTAN EFT - Şube
$ws.Range("C3").Value = ""
$ws.Range("F3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4 - HESAPTAN EFT - ATM
$ws.Range("C4").Value = ""
$ws.Range("F4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5 - HESAPTAN EFT - Mobil
$ws.Range("C5").Value = ""
$ws.Range("F5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 6 - DÜZENLİ EFT
$ws.Range("C6").Value = ""

# Row 7 - KREDİ KARTINDAN FATURA ÖDEME
$ws.Range("F7").Value = "%3"

# Row 8 - HESAPTAN HAVALE - Şube
$ws.Range("C8").Value = ""
$ws.Range("F8").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 9 - HESAPTAN HAVALE - ATM
$ws.Range("C9").Value = ""
$ws.Range("F9").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 10 - HESAPTAN HAVALE - Mobil
$ws.Range("C10").Value = ""
$ws.Range("F10").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 11 - DÜZENLİ HAVALE
$ws.Range("C11").Value = ""

# Row 12 - GİDEN SWIFT
$ws.Range("C12").Value = ""

# Row 13 - GELEN SWIFT
$ws.Range("C13").Value = ""
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"
$ws.Range("F13").Value = "Hesaba: Asgari 300 TL | Azami 3.080 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

# Row 14 - GİDEN SWIFT - Mobil
$ws.Range("C14").Value = ""
$ws.Range("F14").Value = "1.952,38 TL - 9.523,81 TL"
